# T30R09_data.xlsx edit
#
# Summary of the author's change (per commit message / diff):
#  - "1881" survey tab renamed to "1880Survey"
#  - "1942" survey tab renamed to "1940Survey"
#  - Two new metadata tabs added at the end: "1940Metadata" and "1880Metadata",
#    each holding a single explanatory note about why the tab name differs
#    from the actual survey year.
#  - "1940Survey" (old "1942") becomes the active/selected tab instead of
#    "1942notes".

$wb = $excel.ActiveWorkbook

# --- Rename the two data tabs -------------------------------------------------
$surveyTab1880 = $wb.Worksheets.Item(1)   # was "1881"
$surveyTab1880.Name = "1880Survey"

$surveyTab1940 = $wb.Worksheets.Item(3)   # was "1942"
$surveyTab1940.Name = "1940Survey"

# --- Add the two new metadata sheets at the end, in order ---------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)

$metadata1940 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$metadata1940.Name = "1940Metadata"

$metadata1880 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $metadata1940)
$metadata1880.Name = "1880Metadata"

# Fill in the note text. Order matters for shared-string de-duplication, so
# write the 1880/1881 note before the 1940/1942 note.
$metadata1880.Range("A1").Value = "Actually surveyed in 1881; tab reads 1880 for consistency in the data processing R script"
$metadata1880.Range("A1").Font.Color = 0

$metadata1940.Range("A1").Value = "Actually surveyed in 1942; tab reads 1940 for consistency in the data processing R script"
$metadata1940.Range("A1").Font.Color = 0

# Leave the cursor resting on row 2 (below the note) of "1940Metadata".
$metadata1940.Range("A2").Select() | Out-Null

# --- Update selections / active tab -------------------------------------------
# Put the "1880Survey" tab's selection/scroll near the cell the author last
# looked at.
$surveyTab1880.Activate() | Out-Null
$surveyTab1880.Range("I81").Select() | Out-Null

# The "1940Survey" tab (formerly "1942") becomes the active tab, matching the
# workbook's new activeTab / tabSelected state. Its own selection (G2:G76)
# stays as it was.
$surveyTab1940.Activate() | Out-Null
$surveyTab1940.Range("G2:G76").Select() | Out-Null
